# "Actualización 10 de Mayo"
# Updates the statistics on the "Estadisticos 1P", "Estadisticos 2P" and
# "Estadisticos Final" sheets, and populates the "Rescatables" sheet with
# the list of students who need to retake an exam.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Estadisticos 1P
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")

$ws1.Cells.Item(2,4).Value = 0
$ws1.Cells.Item(2,5).Value = 4
$ws1.Cells.Item(2,6).Value = 30
$ws1.Cells.Item(2,7).Value = 88.23999999999999
$ws1.Cells.Item(2,8).Value = 8.1

$ws1.Cells.Item(3,4).Value = 0
$ws1.Cells.Item(3,5).Value = 3
$ws1.Cells.Item(3,8).Value = 8.1

$ws1.Cells.Item(4,4).Value = 0
$ws1.Cells.Item(4,5).Value = 4
$ws1.Cells.Item(4,6).Value = 24
$ws1.Cells.Item(4,7).Value = 85.70999999999999
$ws1.Cells.Item(4,8).Value = 7.9

# ---------------------------------------------------------------------
# Estadisticos 2P
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Cells.Item(2,4).Value = 2
$ws2.Cells.Item(2,5).Value = 5
$ws2.Cells.Item(2,6).Value = 29
$ws2.Cells.Item(2,7).Value = 85.29000000000001
$ws2.Cells.Item(2,8).Value = 8.300000000000001

$ws2.Cells.Item(3,4).Value = 8
$ws2.Cells.Item(3,5).Value = 8
$ws2.Cells.Item(3,6).Value = 20
$ws2.Cells.Item(3,7).Value = 71.43000000000001
$ws2.Cells.Item(3,8).Value = 8.5

$ws2.Cells.Item(4,4).Value = 11
$ws2.Cells.Item(4,5).Value = 11
$ws2.Cells.Item(4,6).Value = 17
$ws2.Cells.Item(4,7).Value = 60.71
$ws2.Cells.Item(4,8).Value = 8.5

# ---------------------------------------------------------------------
# Estadisticos Final
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")

$ws3.Cells.Item(2,4).Value = 0
$ws3.Cells.Item(2,5).Value = 4
$ws3.Cells.Item(2,6).Value = 30
$ws3.Cells.Item(2,7).Value = 88.23999999999999
$ws3.Cells.Item(2,8).Value = 8.4

$ws3.Cells.Item(3,4).Value = 0
$ws3.Cells.Item(3,5).Value = 3
$ws3.Cells.Item(3,8).Value = 8.1

$ws3.Cells.Item(4,4).Value = 0
$ws3.Cells.Item(4,5).Value = 4
$ws3.Cells.Item(4,6).Value = 24
$ws3.Cells.Item(4,7).Value = 85.70999999999999
$ws3.Cells.Item(4,8).Value = 7.9

# ---------------------------------------------------------------------
# Rescatables - list of students pending a make-up exam
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$NC        = @(18330051920446,18330051920446,18330051920190,18330051920190,18330051920391,18330051920085,18330051920103,18330051920189)
$Paterno   = @("CIRUELO","CIRUELO","VERA","VERA","ZACARIAS","ANTONIO","MACUIXTLE","VERA")
$Materno   = @("MANCILLA","MANCILLA","PAZOS","PAZOS","HERNANDEZ","IGNACIO","MACUIXTLE","PAZOS")
$Nombres   = @("MARIA DEL CARMEN","MARIA DEL CARMEN","CARLOS DANIEL","CARLOS DANIEL","LUIS ALBERTO","JOSE MANUEL","JOSUE SAMUEL","CARLA DANIELA")
$NombreLgo = @(
    "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS",
    "DETERMINA LA NÓMINA DEL PERSONAL DE LA ORGANIZACIÓN TOMANDO EN CUENTA LA NORMATIVIDAD LABORAL",
    "DETERMINA LA NÓMINA DEL PERSONAL DE LA ORGANIZACIÓN TOMANDO EN CUENTA LA NORMATIVIDAD LABORAL",
    "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS",
    "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS",
    "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS",
    "DETERMINA LA NÓMINA DEL PERSONAL DE LA ORGANIZACIÓN TOMANDO EN CUENTA LA NORMATIVIDAD LABORAL",
    "DETERMINA REMUNERACIONES DEL PERSONAL EN SITUACIONES EXTRAORDINARIAS"
)
$Grupo     = @("6ARHV","6ARHV","6ARHV","6ARHV","6ARHV","6ARHV","6ARHV","6ARHV")
$Reprob    = @(2,2,2,2,2,1,1,1)

$firstRow = 2

# Column by column so the shared-string table is populated in the same
# order the source workbook used (Paterno, then Materno, then Nombres).
for ($i = 0; $i -lt $Paterno.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 2).Value = $Paterno[$i]
}
for ($i = 0; $i -lt $Materno.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 3).Value = $Materno[$i]
}
for ($i = 0; $i -lt $Nombres.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 4).Value = $Nombres[$i]
}
for ($i = 0; $i -lt $NombreLgo.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 5).Value = $NombreLgo[$i]
}
for ($i = 0; $i -lt $Grupo.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 6).Value = $Grupo[$i]
}
for ($i = 0; $i -lt $NC.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 1).Value = $NC[$i]
}
for ($i = 0; $i -lt $Reprob.Length; $i++) {
    $ws4.Cells.Item($firstRow + $i, 7).Value = $Reprob[$i]
}
